# Update final output. Correcting duplicated guides issue.
#
# The order list contained four groups of guide rows that duplicated
# sequences/IDs already covered by a merged multi-SNP guide group
# (e.g. "rs1317708|rs874424_*" duplicating "rs1317708|rs4989024|rs874424_*",
# "rs2287322_*" and "rs8077638|rs34962442|rs62090051_*" duplicating
# "rs2287322|rs8077638|rs34962442|rs62090051_*", and "rs4989024_*"
# duplicating the same merged group). Remove those 4 redundant
# 3-row (Guide_ID/Sequence _1,_2,_3) blocks entirely.
#
# Delete from the bottom up so earlier row numbers stay valid while
# later rows are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rs8077638|rs34962442|rs62090051_1..3 (duplicate of rs2287322|rs8077638|rs34962442|rs62090051_*)
$ws.Range("A242:B244").EntireRow.Delete()

# rs4989024_1..3 (duplicate of rs1317708|rs4989024|rs874424_*)
$ws.Range("A155:B157").EntireRow.Delete()

# rs2287322_1..3 (duplicate of rs2287322|rs8077638|rs34962442|rs62090051_*)
$ws.Range("A116:B118").EntireRow.Delete()

# rs1317708|rs874424_1..3 (duplicate of rs1317708|rs4989024|rs874424_*)
$ws.Range("A71:B73").EntireRow.Delete()
